# Refresh the cryptocurrency price/volume table with freshly scraped values,
# matching the automated "Updated cryptos list ... with GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.064.39"
$ws.Range("E2").Value = "  -3.37%  "
$ws.Range("D3").Value = "1.847.38"
$ws.Range("E3").Value = "  -2.34%  "
$ws.Range("D4").Value = "'0.9994"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'0.7057"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.15%  "
$ws.Range("D6").Value = "'238.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.96%  "
$ws.Range("D7").Value = "'0.9995"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "'0.3049"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.80%  "
$ws.Range("D9").Value = "'0.07477"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.20%  "
$ws.Range("D10").Value = "'23.38"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.58%  "
$ws.Range("D11").Value = "'0.08128"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.72%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.849.00"
$ws.Range("E12").Value = "  -3.40%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "'0.7259"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.50%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'5.230"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.58%  "
$ws.Range("D15").Value = "'89.18"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.15%  "
$ws.Range("D16").Value = "29.054.77"
$ws.Range("E16").Value = "  -3.65%  "
$ws.Range("D17").Value = "'5.791"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.11%  "
$ws.Range("D18").Value = "'239.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.92%  "
$ws.Range("D19").Value = "'0.000007676"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.38%  "
$ws.Range("D20").Value = "'13.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.33%  "
$ws.Range("D21").Value = "'0.9992"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").Value = "2.083.98"
$ws.Range("E22").Value = "  -3.01%  "
$ws.Range("D23").Value = "'0.9995"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").Value = "'7.555"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.96%  "
$ws.Range("D25").Value = "'0.1464"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.65%  "
$ws.Range("D26").Value = "'8.966"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.67%  "
$ws.Range("D27").Value = "'160.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.01%  "
$ws.Range("E28").Value = "  -4.13%  "
$ws.Range("D29").Value = "'1.938"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.80%  "
$ws.Range("D30").Value = "'1.385"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.13%  "
$ws.Range("D31").Value = "'4.561"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.90%  "
$ws.Range("E32").Value = "  -2.86%  "
$ws.Range("D33").Value = "'4.003"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.59%  "
$ws.Range("D34").Value = "'0.05165"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.05%  "
$ws.Range("D35").Value = "'1.186"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.61%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.7068"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.45%  "
$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").Value = "'1.031"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.98%  "
$ws.Range("D38").Value = "'2.643"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.79%  "
$ws.Range("E39").Value = "  -5.58%  "
$ws.Range("D40").Value = "'2.676"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.29%  "
$ws.Range("D41").Value = "'0.9070"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.12%  "
$ws.Range("D42").Value = "'5.978"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.85%  "
$ws.Range("E43").Value = "  -6.06%  "
$ws.Range("D44").Value = "1.063.81"
$ws.Range("E44").Value = "  -3.39%  "
$ws.Range("D45").Value = "'70.13"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.70%  "
$ws.Range("D46").Value = "'0.9990"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").Value = "'102.14"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.45%  "
$ws.Range("D48").Value = "'1.752"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.31%  "
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").Value = "'7.062"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.41%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'9.178"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.34%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "1.976.05"
$ws.Range("E51").Value = "  -4.28%  "
